# lines_states.xlsx — add line7 / line8 contingency rows and refresh the
# extr1..extr8 rows that follow them ("contingencies with rene fine").
#
# Net effect vs. the original sheet:
#   - rows 8 and 9 change from the first two "extr" rows into new
#     "line7" / "line8" rows with fresh C/D/E data
#   - the data that used to live in rows 8-15 (extr1..extr8) is shifted
#     down by two rows (now rows 10-17) and several C/D/E values are
#     refreshed along the way
#   - two brand-new rows (16, 17) are appended carrying extr7 / extr8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => (A index, B name, C, D, E in_service)
$rows = @{
    8  = @(6,  "line7", 14, 11, $true)
    9  = @(7,  "line8", 16, 9,  $true)
    10 = @(8,  "extr1", 5,  12, $true)
    11 = @(9,  "extr2", 5,  9,  $true)
    12 = @(10, "extr3", 10, 11, $true)
    13 = @(11, "extr4", 7,  8,  $true)
    14 = @(12, "extr5", 9,  11, $true)
    15 = @(13, "extr6", 7,  11, $true)
    16 = @(14, "extr7", 5,  7,  $true)
    17 = @(15, "extr8", 8,  5,  $true)
}

# The two new rows need the same look as the rest of the A column
# (bold, centered, top-aligned, boxed) - just clone the formatting that
# row 15's A cell already carries instead of rebuilding it by hand.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($r in 8..17) {
    $data = $rows[$r]
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]
}
